# Append new daily COVID data rows (44263..44267 -> 2021-03-08 .. 2021-03-12)
# to Planilha1 (sheet1), continuing directly after the existing last row (249).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$data = @(
    @(44263, 6268, 169, 1574, 8011, 1310, 227, 11, 216, 37),
    @(44264, 6306, 186, 1602, 8094, 1329, 236, 12, 224, 37),
    @(44265, 6358, 288, 1629, 8275, 1347, 245, 12, 233, 37),
    @(44266, 6403, 301, 1651, 8355, 1356, 258, 12, 246, 37),
    @(44267, 6467, 276, 1698, 8441, 1377, 284, 14, 270, 37)
)

$startRow = 250
$endRow = $startRow + $data.Count - 1

# Copy formatting (styles / number formats) from the last existing row
# down onto the new rows before writing the values into them.
$ws.Range("A249:J249").Copy()
$ws.Range("A" + $startRow + ":J" + $endRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
